$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The column header in B1 was renamed from "النص المترجم" to "الترجمة"
$ws.Range("B1").Value = "الترجمة"
